# HINDUNILVR.NS.xlsx monthly stock data — add a "backup" column (R) and
# append three new monthly rows (344-346), plus a handful of Q-column
# corrections that were re-detected as part of the same pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New column header: R1 = "backup" (reuse the same bold header style
#    the rest of row 1 already carries).
# ---------------------------------------------------------------------
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("R1").Value = "backup"

# ---------------------------------------------------------------------
# 2. Backfill column R for every existing data row (2-343) with 0,
#    except row 338 which is 1.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 343; $r++) {
    if ($r -eq 338) {
        $ws.Cells.Item($r, 18).Value = 1
    } else {
        $ws.Cells.Item($r, 18).Value = 0
    }
}

# ---------------------------------------------------------------------
# 3. A handful of Q ("detect_structure") values were re-run and reset
#    to 0 on this pass.
# ---------------------------------------------------------------------
$qRowsToZero = @(25, 28, 31, 35, 36, 41, 48, 54)
foreach ($r in $qRowsToZero) {
    $ws.Cells.Item($r, 17).Value = 0
}

# ---------------------------------------------------------------------
# 4. Row 341's O ("isPivot") value flips from 0 to 2.
# ---------------------------------------------------------------------
$ws.Cells.Item(341, 15).Value = 2

# ---------------------------------------------------------------------
# 5. Append three new monthly rows (Jul/Aug/Sep 2024).
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row=344; A=45474; B=2461.050048828125; C=2811.300048828125; D=2450.10009765625;  E=2705.64990234375; F=2705.64990234375; G=48620480; H=2024; I=7; J=1; K=0; L=0; M=0; N=27; O=0; P=0; Q=0 },
    @{ Row=345; A=45505; B=2714;               C=2834.949951171875; D=2666.199951171875; E=2778;              F=2778;              G=36228373; H=2024; I=8; J=1; K=0; L=0; M=0; N=31; O=0; P=0; Q=1 },
    @{ Row=346; A=45536; B=2794;               C=3035;               D=2771.64990234375;  E=2966.25;           F=2966.25;           G=36960781; H=2024; I=9; J=1; K=0; L=0; M=0; N=35; O=0; P=0; Q=0 }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 1).Value  = $rowData.A
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value  = $rowData.B
    $ws.Cells.Item($r, 3).Value  = $rowData.C
    $ws.Cells.Item($r, 4).Value  = $rowData.D
    $ws.Cells.Item($r, 5).Value  = $rowData.E
    $ws.Cells.Item($r, 6).Value  = $rowData.F
    $ws.Cells.Item($r, 7).Value  = $rowData.G
    $ws.Cells.Item($r, 8).Value  = $rowData.H
    $ws.Cells.Item($r, 9).Value  = $rowData.I
    $ws.Cells.Item($r, 10).Value = $rowData.J
    $ws.Cells.Item($r, 11).Value = $rowData.K
    $ws.Cells.Item($r, 12).Value = $rowData.L
    $ws.Cells.Item($r, 13).Value = $rowData.M
    $ws.Cells.Item($r, 14).Value = $rowData.N
    $ws.Cells.Item($r, 15).Value = $rowData.O
    $ws.Cells.Item($r, 16).Value = $rowData.P
    $ws.Cells.Item($r, 17).Value = $rowData.Q
}

Write-Output "done"
